# Update Name of Algo
# Apply updated numeric values to the RandomForest imputation result sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.852199999999995
$ws.Range("A3").Value = -21.42190000000002
$ws.Range("B5").Value = 5.293599999999995
$ws.Range("E5").Value = 13.82029999999999
$ws.Range("E9").Value = 13.21800000000001
$ws.Range("E11").Value = 13.2592
$ws.Range("A14").Value = -20.58499999999998
$ws.Range("A21").Value = -21.38050000000001
$ws.Range("E21").Value = 12.5991
$ws.Range("A23").Value = -21.40070000000003
$ws.Range("A25").Value = -22.44190000000004
